$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.449.57"
$ws.Range("E2").Value = "  +5.01%  "
$ws.Range("D3").Value = "2.495.14"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.40"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.95"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.521"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.52"
$ws.Range("E10").Value = "  +5.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.28"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.13"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "2.886.55"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "2.498.50"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.841"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "47.338.72"
$ws.Range("E18").Value = "  +4.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.54"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").Value = "0.0₃0934"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.73"
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.32"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  +4.99%  "
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.15"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.19"
$ws.Range("E30").Value = "  +6.46%  "
$ws.Range("E31").Value = "  +7.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.46"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.94"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +3.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.23"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.27"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("D45").Value = "1.968.17"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.20"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.33"
$ws.Range("E50").Value = "  +13.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.69"
$ws.Range("E51").Value = "  +2.88%  "
